# DeCuongThucHienDoAn_Python.docx edit script
# Commit: "Bo sung chuc nang gan nhan POS"
#
# Applies, in document order:
#  1. Renumber section 5: the old 5.2 "Thiet ke co so du lieu" paragraph's
#     text is changed so it becomes the new 5.3 "Xay dung chuong trinh";
#     a brand new paragraph is inserted just above it holding the new 5.2
#     "Thiet ke giao dien" (the content that used to live in its own 5.3
#     paragraph); the old 5.3 paragraph is removed; the old 5.4
#     paragraph's text becomes the new 5.4 "Thu nghiem chuong trinh" and
#     the old 5.5 paragraph (now redundant) is removed.
#  2. Collapse four runs of the NLP bibliography entry into two runs.
#  3. Mark the "Da Lat ngay ..." run with a lastRenderedPageBreak.
#
# Work from the bottom of the document upwards so paragraph indices for
# not-yet-processed items stay valid.

$d = $word.ActiveDocument

function Fill-ParagraphXml($paragraph, [string]$bodyFragment) {
    $xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
$bodyFragment
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
    $paragraph.Range.InsertXML($xml)
}

function Fill-RangeXml($range, [string]$bodyFragment) {
    $xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
$bodyFragment
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
    $range.InsertXML($xml)
}

# ---------------------------------------------------------------------
# Locate the relevant paragraphs by distinctive text before any edits.
# ---------------------------------------------------------------------
$idxDaLat = 0
$idxNlp = 0
$idx52 = 0
$idx53 = 0
$idx54 = 0
$idx55 = 0

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -match "L.t ng.y 7 th.ng 3") { $idxDaLat = $i }
    if ($t -match "ng.n ng. t. nhi.n") { $idxNlp = $i }
    if ($t -match "5\.2\. Thi.t k. c. s. d. li.u") { $idx52 = $i }
    if ($t -match "^5\.3\. Thi.t k. giao di.n") { $idx53 = $i }
    if ($t -match "5\.4\. X.y d.ng ch..ng tr.nh") { $idx54 = $i }
    if ($t -match "5\.5\. Th. nghi.m ch..ng tr.nh") { $idx55 = $i }
}

Write-Host "Found paragraphs:" $idx52 $idx53 $idx54 $idx55 $idxNlp $idxDaLat

# ---------------------------------------------------------------------
# 3) "Da Lat ngay 7 thang 3 nam 2022" -- add lastRenderedPageBreak
# ---------------------------------------------------------------------
$pDaLat = $d.Paragraphs.Item($idxDaLat)
Fill-ParagraphXml $pDaLat @"
<w:p><w:pPr><w:jc w:val="right"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:i/><w:iCs/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:i/><w:iCs/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:lastRenderedPageBreak/><w:t>&#272;&#224; L&#7841;t ng&#224;y 7 th&#225;ng 3 n&#259;m 2022</w:t></w:r></w:p>
"@

# ---------------------------------------------------------------------
# 2) NLP bibliography entry -- merge the ":"+" "+title+authors runs into
#    two runs, leaving the leading "[6] Ebook ve ngon ngu tu nhien" run
#    and the paragraph's pPr untouched.
# ---------------------------------------------------------------------
$pNlp = $d.Paragraphs.Item($idxNlp)
$rNlp = $pNlp.Range
$nlpTarget = $d.Range($rNlp.Start + 30, $rNlp.End)
Fill-RangeXml $nlpTarget @"
<w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="202124"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Natural Language Processing with Python &#8211; Steven Bird, Ewan Klein, Edward Loper</w:t></w:r></w:p>
"@

# ---------------------------------------------------------------------
# 1) Section 5 renumbering
# ---------------------------------------------------------------------

# old 5.5 "Thu nghiem chuong trinh" paragraph is dropped entirely -- its
# text is absorbed into the renumbered old-5.4 paragraph below.
$d.Paragraphs.Item($idx55).Range.Delete()

# old 5.4 "Xay dung chuong trinh" paragraph keeps its tab but becomes the
# new 5.4 "Thu nghiem chuong trinh".
$p54 = $d.Paragraphs.Item($idx54)
Fill-ParagraphXml $p54 @"
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="1080"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:tab/><w:t>5.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>4</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>. Th&#7917; nghi&#7879;m ch&#432;&#417;ng tr&#236;nh</w:t></w:r></w:p>
"@

# old 5.3 "Thiet ke giao dien" paragraph is removed entirely -- its text is
# reused (renumbered) in the brand-new paragraph inserted above old 5.2.
$d.Paragraphs.Item($idx53).Range.Delete()

# old 5.2 "Thiet ke co so du lieu" paragraph keeps its leading tab + "5."
# run but becomes the new 5.3 "Xay dung chuong trinh ".
$p52 = $d.Paragraphs.Item($idx52)
Fill-ParagraphXml $p52 @"
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="1080"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:tab/><w:t>5.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>3</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>. X&#226;y d&#7921;ng ch&#432;&#417;ng tr&#236;nh</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p>
"@

# Insert a brand new, empty paragraph right before the (just renumbered)
# old-5.2 paragraph, then fill it in with the new 5.2 "Thiet ke giao dien"
# content.
$p52.Range.InsertParagraphBefore()
$pNewIdx = $idx52
$pNew = $d.Paragraphs.Item($pNewIdx)
Fill-ParagraphXml $pNew @"
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="1080" w:firstLine="360"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>5.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>2</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">. </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Thi&#7871;t k&#7871; giao di&#7879;n</w:t></w:r></w:p>
"@

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
